# Update cryptocurrency price/volume data per Dec 23 2023 refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.670.19'
$ws.Range("E2").Value = '  -0.66%  '

$ws.Range("D3").Value = '2.277.27'
$ws.Range("E3").Value = '  +1.35%  '

$ws.Range("E4").Value = '  +0.18%  '

$ws.Range("D5").Value = "'94.91"
$ws.Range("E5").Value = '  -2.72%  '

$ws.Range("D6").Value = "'267.13"
$ws.Range("E6").Value = '  -2.01%  '

$ws.Range("E7").Value = '  +0.56%  '

$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("E9").Value = '  -3.80%  '

$ws.Range("D10").Value = "'44.22"
$ws.Range("E10").Value = '  -7.59%  '

$ws.Range("E11").Value = '  -0.81%  '

$ws.Range("D12").Value = "'7.72"
$ws.Range("E12").Value = '  -7.39%  '

$ws.Range("D13").Value = "'0.104"
$ws.Range("E13").Value = '  -0.26%  '

$ws.Range("D14").Value = '2.619.67'
$ws.Range("E14").Value = '  +1.47%  '

$ws.Range("D15").Value = "'15.13"
$ws.Range("E15").Value = '  -1.24%  '

$ws.Range("E16").Value = '  +1.88%  '

$ws.Range("D17").Value = '2.284.13'
$ws.Range("E17").Value = '  +1.92%  '

$ws.Range("D18").Value = '43.587.16'
$ws.Range("E18").Value = '  -0.85%  '

$ws.Range("E19").Value = '  +1.25%  '

$ws.Range("D20").Value = "'6.15"
$ws.Range("E20").Value = '  -0.97%  '

$ws.Range("D21").Value = "'71.97"
$ws.Range("E21").Value = '  +1.47%  '

$ws.Range("D22").Value = "'2.35"
$ws.Range("E22").Value = '  +0.39%  '

$ws.Range("D23").Value = "'234.53"
$ws.Range("E23").Value = '  +0.00%  '

$ws.Range("D24").Value = "'8.91"
$ws.Range("E24").Value = '  -5.51%  '

$ws.Range("E25").Value = '  -0.05%  '

$ws.Range("D26").Value = "'11.33"
$ws.Range("E26").Value = '  -1.01%  '

$ws.Range("E27").Value = '  -0.74%  '

$ws.Range("E28").Value = '  -2.25%  '

$ws.Range("E29").Value = '  -0.05%  '

$ws.Range("D30").Value = "'39.09"
$ws.Range("E30").Value = '  -1.50%  '

$ws.Range("D31").Value = "'176.47"
$ws.Range("E31").Value = '  +1.75%  '

$ws.Range("E32").Value = '  +3.77%  '

$ws.Range("D33").Value = "'0.0880"
$ws.Range("E33").Value = '  -3.51%  '

$ws.Range("D34").Value = "'5.34"
$ws.Range("E34").Value = '  -5.04%  '

$ws.Range("E35").Value = '  +0.49%  '

$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").Value = "'0.107"
$ws.Range("E36").Value = '  -4.50%  '

$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = "'0.0354"
$ws.Range("E37").Value = '  +0.75%  '

$ws.Range("D38").Value = "'4.41"
$ws.Range("E38").Value = '  +0.60%  '

$ws.Range("E39").Value = '  -9.68%  '

$ws.Range("D40").Value = "'2.35"
$ws.Range("E40").Value = '  +7.54%  '

$ws.Range("E41").Value = '  -6.29%  '

$ws.Range("D42").Value = "'1.34"
$ws.Range("E42").Value = '  +16.06%  '

$ws.Range("D43").Value = "'11.81"
$ws.Range("E43").Value = '  -5.81%  '

$ws.Range("D44").Value = "'62.15"
$ws.Range("E44").Value = '  +0.00%  '

$ws.Range("D45").Value = "'8.78"
$ws.Range("E45").Value = '  +4.16%  '

$ws.Range("D46").Value = "'5.20"
$ws.Range("E46").Value = '  -4.46%  '

$ws.Range("E47").Value = '  -0.61%  '

$ws.Range("D48").Value = "'98.15"
$ws.Range("E48").Value = '  -2.56%  '

$ws.Range("D49").Value = "'1.18"
$ws.Range("E49").Value = '  -0.67%  '

$ws.Range("E50").Value = '  +6.02%  '

$ws.Range("D51").Value = '2.499.25'
$ws.Range("E51").Value = '  +1.57%  '
